$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the previous row (row 9) down into the new row 10
# so the new cells reuse the existing style indices instead of Excel
# minting new ones.
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new log entry
$ws.Cells.Item(10, 1).Value = 43131
$ws.Cells.Item(10, 2).Value = "Drew"
$ws.Cells.Item(10, 3).Value = 30
$ws.Cells.Item(10, 4).Value = "Programming"
$ws.Cells.Item(10, 5).Value = "Added RenderingSystem and completed Basic ECS"

# Move the selection to the new last cell, matching the saved view state
$ws.Range("E10").Select()
